$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "41.340.77"
$ws.Cells.Item(2, 5).Value = "  -3.30%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.471.44"
$ws.Cells.Item(3, 5).Value = "  -2.60%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.07%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "311.95"
$ws.Cells.Item(5, 5).Value = "  +0.26%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "94.61"
$ws.Cells.Item(6, 5).Value = "  -6.00%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -3.31%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.07%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -4.31%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "33.52"
$ws.Cells.Item(10, 5).Value = "  -6.09%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.78%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.52%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -4.41%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.851.11"
$ws.Cells.Item(14, 5).Value = "  -2.66%  "

# Row 15
$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.06"
$ws.Cells.Item(15, 5).Value = "  -1.87%  "

# Row 16
$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "2.447.59"
$ws.Cells.Item(16, 5).Value = "  -3.09%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.791"
$ws.Cells.Item(17, 5).Value = "  -3.05%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "41.346.16"
$ws.Cells.Item(18, 5).Value = "  -3.28%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -6.05%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.0₃0925"
$ws.Cells.Item(20, 5).Value = "  -2.84%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.25"
$ws.Cells.Item(21, 5).Value = "  -8.80%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "68.53"
$ws.Cells.Item(22, 5).Value = "  -2.28%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "237.28"

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.75"
$ws.Cells.Item(24, 5).Value = "  -4.68%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.09%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -6.39%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "24.04"
$ws.Cells.Item(27, 5).Value = "  -5.62%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -4.32%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.68"
$ws.Cells.Item(29, 5).Value = "  -4.90%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "36.60"
$ws.Cells.Item(30, 5).Value = "  -5.24%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "151.72"
$ws.Cells.Item(31, 5).Value = "  -4.62%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "5.48"

# Row 33
$ws.Cells.Item(33, 2).Value = "WEMIXToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.59"
$ws.Cells.Item(33, 5).Value = "  -3.05%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "ApeXProtocol"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.58"
$ws.Cells.Item(34, 5).Value = "  -6.52%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0747"
$ws.Cells.Item(35, 5).Value = "  -5.76%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.07"
$ws.Cells.Item(36, 5).Value = "  -2.46%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "17.48"
$ws.Cells.Item(37, 5).Value = "  -3.49%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -4.86%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -2.96%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Kaspa"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.102"
$ws.Cells.Item(40, 5).Value = "  -8.49%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "4.26"
$ws.Cells.Item(41, 5).Value = "  +3.16%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +0.11%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "19.66"
$ws.Cells.Item(43, 5).Value = "  -9.98%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "1.980.61"
$ws.Cells.Item(44, 5).Value = "  -0.71%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0286"
$ws.Cells.Item(45, 5).Value = "  -4.42%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.02"
$ws.Cells.Item(46, 5).Value = "  -8.40%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -5.01%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "2.714.14"
$ws.Cells.Item(48, 5).Value = "  -2.49%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "69.69"
$ws.Cells.Item(49, 5).Value = "  -3.66%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "97.02"
$ws.Cells.Item(50, 5).Value = "  -4.02%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "74.57"
$ws.Cells.Item(51, 5).Value = "  -6.76%  "
